# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45172 (2023-09-03) to 45175 (2023-09-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$range = $ws.Range("C2:C$lastRow")
$range.Value2 = 45175
